$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.493.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.680.63'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5325'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.28%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +4.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06412'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.42'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07785'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.692.89'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.505'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5621'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅8387'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.85'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.534.23'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.805'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.73'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.388'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.02%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.50'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1274'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.465'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.81%  '
$ws.Range('E28').Value = '  +2.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06132'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.279'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.607'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.458'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.704'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.02%  '
$ws.Range('E34').Value = '  +4.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.792'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.38%  '
$ws.Range('E36').Value = '  +1.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5703'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01639'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.952'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8666'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.058.55'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.88'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.831.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.20%  '
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('E46').Value = '  +5.66%  '
$ws.Range('E47').Value = '  +2.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9993'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05207'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.062'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.69%  '
$ws.Range('E51').Value = '  +0.26%  '
